# Add the new "Dying Light: The Beast PS5" product listing and remove the
# now-obsolete "Dying Light The Beast" placeholder row from the upcoming
# games section.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the placeholder row for "Dying Light The Beast" (row 486) - the
# game has now shipped and gets a full listing instead, further down in
# the sheet. Deleting the row shifts everything below it up by one.
$ws.Rows.Item(486).Delete()

# Append the new product row at the end of the sheet (now row 525) with
# the game title in column A and its product-page link text in column E,
# matching the format used by the other entries in this tail section.
$lastRow = $ws.UsedRange.Rows.Count + 1
$ws.Cells.Item($lastRow, 5).Value2 = "https://arenapsgm.ru/tproduct/555739873852-dying-light-the-beast-ps5"
$ws.Cells.Item($lastRow, 1).Value2 = "Dying Light: The Beast PS5"

# Reflect the final cursor/selection state left behind by the edit: the
# freshly added last row is selected, with the view scrolled back up a
# little from where it had been.
$selRange = "A" + $lastRow + ":XFD" + $lastRow
$ws.Range($selRange).Select()
